{"js": "// Load the last paragraph of the document body (the one that ends with\n// \"A sz\u00e1ll\u00edt\u00e1si \u00e9s sz\u00e1ml\u00e1z\u00e1si adatokat k\u00fcl\u00f6n le kell t\u00e1rolni \").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert a new empty paragraph right after it.\n// NOTE: spaceAfter is expressed in points (same as w:spacing/@w:after divided\n// by 20, since w:after is stored in twentieths of a point). 80 / 20 = 4pt.\nconst emptyParagraph = lastParagraph.insertParagraph(\"\", \"After\");\nemptyParagraph.spaceAfter = 4;\n\n// Insert another new paragraph after the empty one with the new sentence.\nconst textParagraph = emptyParagraph.insertParagraph(\n  \"Bele kell venni, hogy a kiszerel\u00e9seknek is kellenek k\u00fcl\u00f6n dolgok. Bev\u00e9telez\u00e9s belev\u00e9tele.\",\n  \"After\"\n);\ntextParagraph.spaceAfter = 4;\n\nawait context.sync();\n", "ps1": "# Locate the last paragraph of the document body (the one ending with\n# \"A sz\u00e1ll\u00edt\u00e1si \u00e9s sz\u00e1ml\u00e1z\u00e1si adatokat k\u00fcl\u00f6n le kell t\u00e1rolni \").\n$d = $word.ActiveDocument\n$lastPara = $d.Paragraphs.Last\n$r = $lastPara.Range\n$r.Collapse(0)  # wdCollapseEnd\n\n# Insert a new, empty paragraph right after it.\n$r.InsertParagraphAfter()\n$emptyPara = $d.Paragraphs.Last\n$emptyPara.SpaceAfter = 4  # w:spacing w:after=\"80\" (80 twentieths-of-a-point = 4pt)\n\n# Insert another new paragraph after the empty one, with the new sentence.\n$r2 = $emptyPara.Range\n$r2.Collapse(0)\n$r2.InsertParagraphAfter()\n$textPara = $d.Paragraphs.Last\n$textPara.Range.InsertAfter(\"Bele kell venni, hogy a kiszerel\u00e9seknek is kellenek k\u00fcl\u00f6n dolgok. Bev\u00e9telez\u00e9s belev\u00e9tele.\")\n$textPara.SpaceAfter = 4\n"}
